$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 5
$ws.Range("H5").Value = 220
$ws.Range("I5").Value = 220
$ws.Range("K5").Value = 220
$ws.Range("M5").Value = -105

# ALC row 11
$ws.Range("H11").Value = 41668190
$ws.Range("I11").Value = 41668190
$ws.Range("K11").Value = 41668190
$ws.Range("M11").Value = -41668050

# ALC row 43
$ws.Range("H43").Value = 7721.7144
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 8508.666999999999
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 8508.666999999999
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -8646.666999999999

# ALC row 58
$ws.Range("H58").Value = 2545.6
$ws.Range("I58").Value = 465
$ws.Range("J58").Value = 3932.6667
$ws.Range("K58").Value = 1395
$ws.Range("L58").Value = 11798.0001
$ws.Range("M58").Value = -1245
$ws.Range("N58").Value = -12098.0001

# ALC row 86
$ws.Range("H86").Value = 9127.762000000001
$ws.Range("I86").Value = 9543.125
$ws.Range("J86").Value = 7798.6
$ws.Range("K86").Value = 9543.125
$ws.Range("L86").Value = 7798.6
$ws.Range("M86").Value = -8420.125
$ws.Range("N86").Value = -10044.6

# ALC row 89
$ws.Range("H89").Value = 9127.762000000001
$ws.Range("I89").Value = 9543.125
$ws.Range("J89").Value = 7798.6
$ws.Range("K89").Value = 47715.625
$ws.Range("L89").Value = 38993
$ws.Range("M89").Value = -42099.625
$ws.Range("N89").Value = -50225

# ALC row 116
$ws.Range("H116").Value = 6642.533
$ws.Range("I116").Value = 5861.1113
$ws.Range("K116").Value = 5861.1113
$ws.Range("M116").Value = -2419.1113

# ALC row 132
$ws.Range("H132").Value = 2178.5854
$ws.Range("I132").Value = 1966.2778
$ws.Range("J132").Value = 3707.2
$ws.Range("K132").Value = 5898.8334
$ws.Range("L132").Value = 11121.6
$ws.Range("M132").Value = -3368.8334
$ws.Range("N132").Value = -16181.6

# ALC row 137
$ws.Range("H137").Value = 6399.9165
$ws.Range("I137").Value = 1471.75
$ws.Range("K137").Value = 4415.25
$ws.Range("M137").Value = -1865.25

# ALC row 138
$ws.Range("H138").Value = 2201.95
$ws.Range("J138").Value = 3471.0476
$ws.Range("L138").Value = 10413.1428
$ws.Range("N138").Value = -20693.1428

$ws = $wb.Worksheets.Item("ARM")
# ARM row 4
$ws.Range("H4").Value = 1638.6
$ws.Range("I4").Value = 600
$ws.Range("K4").Value = 600
$ws.Range("M4").Value = -484

# ARM row 5
$ws.Range("H5").Value = 126.875
$ws.Range("I5").Value = 102.14286
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 102.14286
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = 9.857140000000001
$ws.Range("N5").Value = -524

# ARM row 32
$ws.Range("H32").Value = 12196741
$ws.Range("I32").Value = 12196741
$ws.Range("K32").Value = 12196741
$ws.Range("M32").Value = -12196454

# ARM row 61
$ws.Range("H61").Value = 26376620
$ws.Range("I61").Value = 41671868
$ws.Range("J61").Value = 156192.28
$ws.Range("K61").Value = 41671868
$ws.Range("L61").Value = 156192.28
$ws.Range("M61").Value = -41671656
$ws.Range("N61").Value = -156616.28

# ARM row 74
$ws.Range("H74").Value = 11372623
$ws.Range("J74").Value = 20773.111
$ws.Range("L74").Value = 20773.111
$ws.Range("N74").Value = -22521.111

# ARM row 77
$ws.Range("H77").Value = 11372623
$ws.Range("J77").Value = 20773.111
$ws.Range("L77").Value = 103865.555
$ws.Range("N77").Value = -112601.555

# ARM row 136
$ws.Range("H136").Value = 26376620
$ws.Range("I136").Value = 41671868
$ws.Range("J136").Value = 156192.28
$ws.Range("K136").Value = 125015604
$ws.Range("L136").Value = 468576.84
$ws.Range("M136").Value = -125013054
$ws.Range("N136").Value = -473676.84

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4
$ws.Range("H4").Value = 126.875
$ws.Range("I4").Value = 102.14286
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 102.14286
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = 12.85714
$ws.Range("N4").Value = -530

# BSM row 99
$ws.Range("H99").Value = 3074.1875
$ws.Range("I99").Value = 2483.8
$ws.Range("J99").Value = 4058.1667
$ws.Range("K99").Value = 2483.8
$ws.Range("L99").Value = 4058.1667
$ws.Range("M99").Value = -985.8000000000002
$ws.Range("N99").Value = -7054.1667

# BSM row 134
$ws.Range("H134").Value = 64446.883
$ws.Range("J134").Value = 152558.14
$ws.Range("L134").Value = 457674.42
$ws.Range("N134").Value = -462744.42

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws.Range("H7").Value = 2436
$ws.Range("I7").Value = 57.272728
$ws.Range("K7").Value = 57.272728
$ws.Range("M7").Value = 55.727272

# CRP row 17
$ws.Range("H17").Value = 28342
$ws.Range("I17").Value = 25008
$ws.Range("J17").Value = 30009
$ws.Range("K17").Value = 25008
$ws.Range("L17").Value = 30009
$ws.Range("M17").Value = -24834
$ws.Range("N17").Value = -30357

# CRP row 22
$ws.Range("H22").Value = 846.82355
$ws.Range("I22").Value = 884.1875
$ws.Range("J22").Value = 249
$ws.Range("K22").Value = 884.1875
$ws.Range("L22").Value = 249
$ws.Range("M22").Value = -534.1875
$ws.Range("N22").Value = -949

# CRP row 33
$ws.Range("H33").Value = 4415.5
$ws.Range("I33").Value = 4415.5
$ws.Range("K33").Value = 4415.5
$ws.Range("M33").Value = -4036.5

# CRP row 99
$ws.Range("H99").Value = 2366.4443
$ws.Range("I99").Value = 2037.375
$ws.Range("K99").Value = 2037.375
$ws.Range("M99").Value = -539.375

# CRP row 105
$ws.Range("H105").Value = 2007.8667
$ws.Range("I105").Value = 1576.1428
$ws.Range("J105").Value = 2385.625
$ws.Range("K105").Value = 1576.1428
$ws.Range("L105").Value = 2385.625
$ws.Range("M105").Value = 170.8571999999999
$ws.Range("N105").Value = -5879.625

# CRP row 126
$ws.Range("H126").Value = 2366.4443
$ws.Range("I126").Value = 2037.375
$ws.Range("K126").Value = 6112.125
$ws.Range("M126").Value = -3642.125

$ws = $wb.Worksheets.Item("CUL")
# CUL row 26
$ws.Range("H26").Value = 169
$ws.Range("J26").Value = 183.33333
$ws.Range("L26").Value = 549.99999
$ws.Range("N26").Value = -1125.99999

# CUL row 51
$ws.Range("H51").Value = 13488
$ws.Range("I51").Value = 8438.444
$ws.Range("K51").Value = 25315.332
$ws.Range("M51").Value = -24855.332

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132
$ws.Range("H132").Value = 55559120
$ws.Range("I132").Value = 66670388
$ws.Range("J132").Value = 2796.3333
$ws.Range("K132").Value = 200011164
$ws.Range("L132").Value = 8388.999899999999
$ws.Range("M132").Value = -200008634
$ws.Range("N132").Value = -13448.9999

# GSM row 134
$ws.Range("H134").Value = 42236.5
$ws.Range("J134").Value = 42236.5
$ws.Range("L134").Value = 126709.5
$ws.Range("N134").Value = -131779.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 1788.7916
$ws.Range("I16").Value = 1581.3077
$ws.Range("J16").Value = 2034
$ws.Range("K16").Value = 1581.3077
$ws.Range("L16").Value = 2034
$ws.Range("M16").Value = -1411.3077
$ws.Range("N16").Value = -2374

# LTW row 93
$ws.Range("H93").Value = 83335130
$ws.Range("I93").Value = 111112170
$ws.Range("K93").Value = 111112170
$ws.Range("M93").Value = -111110922

# LTW row 132
$ws.Range("H132").Value = 508882.03
$ws.Range("I132").Value = 627093.9
$ws.Range("J132").Value = 193650.5
$ws.Range("K132").Value = 1881281.7
$ws.Range("L132").Value = 580951.5
$ws.Range("M132").Value = -1878751.7
$ws.Range("N132").Value = -586011.5

# LTW row 140
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

$ws = $wb.Worksheets.Item("WVR")
# WVR row 4
$ws.Range("H4").Value = 46208.5
$ws.Range("I4").Value = 116000
$ws.Range("J4").Value = 11312.75
$ws.Range("K4").Value = 116000
$ws.Range("L4").Value = 11312.75
$ws.Range("M4").Value = -115887
$ws.Range("N4").Value = -11538.75

# WVR row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

# WVR row 74
$ws.Range("H74").Value = 21666.666
$ws.Range("J74").Value = 21666.666
$ws.Range("L74").Value = 21666.666
$ws.Range("N74").Value = -23538.666

# WVR row 77
$ws.Range("H77").Value = 21666.666
$ws.Range("J77").Value = 21666.666
$ws.Range("L77").Value = 64999.99800000001
$ws.Range("N77").Value = -74359.99800000001

# WVR row 93
$ws.Range("H93").Value = 81954.664
$ws.Range("J93").Value = 77932
$ws.Range("L93").Value = 77932
$ws.Range("N93").Value = -82924

# WVR row 132
$ws.Range("H132").Value = 11729.25
$ws.Range("I132").Value = 2161.5625
$ws.Range("K132").Value = 6484.6875
$ws.Range("M132").Value = -3954.6875

# WVR row 136
$ws.Range("H136").Value = 18970.643
$ws.Range("I136").Value = 1777.7142
$ws.Range("J136").Value = 36163.57
$ws.Range("K136").Value = 5333.142599999999
$ws.Range("L136").Value = 108490.71
$ws.Range("M136").Value = -2783.142599999999
$ws.Range("N136").Value = -113590.71
